$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.595.21"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "2.288.98"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0961"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("E11").Value = "  +4.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.06%  "
$ws.Range("D14").Value = "2.692.81"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "54.577.78"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").Value = "2.308.09"
$ws.Range("E17").Value = "  +3.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.82%  "
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "306.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").Value = "0.0₃0712"
$ws.Range("E28").Value = "  +5.90%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  +4.57%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.936"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.89%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.06%  "
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "125.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0494"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "248.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0898"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.549"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.375"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("E51").Value = "  +7.28%  "
